# Generate Report for Handback
# Updates the "Correspond Handback DateTime" (column H) on the es-es sheet
# for every row whose handback run timestamp was "2016-03-24 18:03:26" (or,
# for a handful of rows that already shared the "2016-03-24 20:14:11" slot,
# were re-handed-back together) to the new handback run timestamp
# "2016-03-30 17:59:27". Also refreshes the "Correspond Handoff Datetime"
# for the TOC.md row (row 19) to its latest handoff time.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("es-es")

$newHandbackTime = "2016-03-30 17:59:27"

$rows = @(2,3,4,7,8,9,10,11,12,13,14,16,18,19,24,25,26,27,28,29,30,31,32,33,34,35,36,41,42,43,44,45,46,47,48,49,50,51,52,53,54,55,56,57,58,59,60,61,62,63)

foreach ($r in $rows) {
    $ws.Range("H" + $r).Value = $newHandbackTime
}

$ws.Range("E19").Value = "2016-03-07 19:32:30"
